$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 for "dct:creator" / "Joseph"
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = "dct:creator"
$ws.Cells.Item(14, 2).Value = "Joseph"

# Row 23 (formerly an empty "vars:" placeholder row, now shifted down by
# the insert above) becomes the new "vars:datascientist" entry. No further
# row insertion is needed here because the template rows 22-90 were all
# identical blank "vars:" placeholders.
$ws.Cells.Item(23, 1).Value = "vars:datascientist"
$ws.Cells.Item(23, 2).Value = "datascientist"
$ws.Cells.Item(23, 5).Value = "a person that knows stuff about data"
